$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Despezas")

# Update existing row 2
$ws.Range("A2").Value = "R$ 2000.00"
$ws.Range("B2").Value = "Trabalho"
$ws.Range("C2").Value = "Ganho"

# Add new row 3
$ws.Range("A3").Value = "R$ 1500"
$ws.Range("B3").Value = "Trabalho 2"
$ws.Range("C3").Value = "Ganho"

# Add new row 4
$ws.Range("A4").Value = "R$ 560.66"
$ws.Range("B4").Value = "Mercado"
$ws.Range("C4").Value = "Gasto"

# Add new row 5
$ws.Range("A5").Value = "R$ 145.60"
$ws.Range("B5").Value = "Carro"
$ws.Range("C5").Value = "Gasto"
